# Updated symbol list on Wed Jan 25 13:57:52 UTC 2023 with GitHub Actions
#
# Refresh the cryptos price table: Price/Volume(1h) are text cells that
# happen to look numeric, so a leading "'" forces Excel to keep them as
# literal text (matching how the sheet was originally authored) instead of
# auto-converting them to numbers. Coin/Link (columns B/C) are plain text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.98"
$ws.Range("E2").Value = "'-3.81%"
$ws.Range("D3").Value = "'35.49"
$ws.Range("E3").Value = "'1.44%"
$ws.Range("D4").Value = "'5.052"
$ws.Range("E4").Value = "'-1.35%"
$ws.Range("D5").Value = "'0.08007"
$ws.Range("E5").Value = "'-1.71%"
$ws.Range("E6").Value = "'-9.04%"
$ws.Range("D7").Value = "'7.811"
$ws.Range("E7").Value = "'-1.75%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'4.058"
$ws.Range("E8").Value = "'-2.13%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.958"
$ws.Range("E9").Value = "'2.10%"
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").Value = "'0.9254"
$ws.Range("E10").Value = "'-0.78%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.1342"
$ws.Range("E11").Value = "'30.23%"
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1899"
$ws.Range("E12").Value = "'-1.55%"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.09225"
$ws.Range("E13").Value = "'1.63%"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03384"
$ws.Range("E14").Value = "'-7.43%"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09867"
$ws.Range("E15").Value = "'-0.21%"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001395"
$ws.Range("E16").Value = "'-2.86%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005762"
$ws.Range("E17").Value = "'-1.54%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.511"
$ws.Range("E18").Value = "'1.24%"
$ws.Range("E19").Value = "'-0.17%"
$ws.Range("D20").Value = "'0.1303"
$ws.Range("E20").Value = "'-2.16%"
$ws.Range("D21").Value = "'5.068"
$ws.Range("E21").Value = "'-0.96%"
$ws.Range("D23").Value = "'0.04495"
$ws.Range("E23").Value = "'-1.24%"
$ws.Range("D24").Value = "'0.001216"
$ws.Range("E24").Value = "'-2.67%"
$ws.Range("D25").Value = "'0.004798"
$ws.Range("E25").Value = "'2.16%"
$ws.Range("D26").Value = "'0.0001252"
$ws.Range("E26").Value = "'-0.14%"
$ws.Range("D27").Value = "'0.0003005"
$ws.Range("E27").Value = "'-33.34%"
$ws.Range("D39").Value = "'0.01908"
$ws.Range("E39").Value = "'-2.05%"
$ws.Range("D40").Value = "'0.04735"
$ws.Range("E40").Value = "'-3.34%"
$ws.Range("D41").Value = "'0.007371"
$ws.Range("E41").Value = "'-3.24%"
$ws.Range("D42").Value = "'0.009675"
$ws.Range("E42").Value = "'22.70%"
$ws.Range("D43").Value = "'0.1332"
$ws.Range("E43").Value = "'-4.04%"
$ws.Range("D44").Value = "'0.002113"
$ws.Range("D45").Value = "'0.01053"
$ws.Range("E45").Value = "'-10.67%"
$ws.Range("D46").Value = "'0.00006348"
$ws.Range("E46").Value = "'-6.05%"
$ws.Range("E47").Value = "'-0.09%"
$ws.Range("E48").Value = "'-65.45%"
$ws.Range("D49").Value = "'0.001662"
$ws.Range("E49").Value = "'-2.44%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'-0.09%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'-0.09%"
